# FALanguage新增标签.xlsx — add 13 new trigger-editor / compact-mode label rows
# (can open up to 10 trigger editor; can drag input params)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Key / Chinese / English rows appended right after the existing data
# (existing data ends at row 838), starting at row 839.
$data = @(
    @('TriggerOpenNewEditor', '多开编辑器', 'New editor'),
    @('TriggerNewShort', '新', 'N'),
    @('TriggerPlaceOnMapShort', '放置', 'PM'),
    @('TriggerDisabledShort', '禁', 'D'),
    @('TriggerEasyShort', '简', 'E'),
    @('TriggerMediumShort', '普', 'N'),
    @('TriggerHardShort', '困', 'H'),
    @('TriggerAddShort', '新', 'N'),
    @('TriggerCloneShort', '复', 'C'),
    @('TriggerDeleteShort', '删', 'D'),
    @('SearchReferenceTitleShort', '查找', 'SR'),
    @('TriggerCompactMode', '紧凑模式', 'Compact Mode'),
    @('TriggerCompactModeShort', '紧凑', 'Compact'),
)

$startRow = 839
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $triple = $data[$i]
    if ($row -eq 850) {
        # This particular row was authored with the Chinese label typed
        # before the English key (matches original shared-string order).
        $ws.Cells.Item($row, 2).Value = $triple[1]
        $ws.Cells.Item($row, 1).Value = $triple[0]
        $ws.Cells.Item($row, 3).Value = $triple[2]
    } else {
        $ws.Cells.Item($row, 1).Value = $triple[0]
        $ws.Cells.Item($row, 2).Value = $triple[1]
        $ws.Cells.Item($row, 3).Value = $triple[2]
    }
}

# Scroll the view down near the bottom of the sheet and select the newly
# added last row, matching where the author was working when they saved.
$win = $excel.ActiveWindow
$win.ScrollRow = 833
$win.ScrollColumn = 1

$ws.Range("A850:C851").Select()
